$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("bot")
$ws9 = $wb.Worksheets.Item("9")
$ws1.Move($null, $ws9)
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
